$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Give them the same format as the other header cells (A1:F1) by
# copying the format from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Tiny precision corrections to the existing MSE / MAE values
$ws.Range("B2").Value = 0.4108830721519875
$ws.Range("D2").Value = 0.525669270734418

# New data values for the added columns
$ws.Range("G2").Value = 0.1228586025167412
$ws.Range("H2").Value = 0.991
